$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.810.18'
$ws.Range("E2").Value = '  +1.49%  '

$ws.Range("D3").Value = '3.009.97'
$ws.Range("E3").Value = '  +3.56%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '383.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.28%  '

$ws.Range("E7").Value = '  +1.34%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.600'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.61'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.19%  '

$ws.Range("E11").Value = '  +0.29%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0850'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.84%  '

$ws.Range("D13").Value = '3.473.51'
$ws.Range("E13").Value = '  +3.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.54'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.54%  '

$ws.Range("E15").Value = '  +3.10%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.02'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +11.06%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.993.47'
$ws.Range("E17").Value = '  +3.01%  '

$ws.Range("D18").Value = '51.761.85'
$ws.Range("E18").Value = '  +1.57%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.32'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.99%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.49'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.03'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.04%  '

$ws.Range("D22").Value = '0.0₃0968'
$ws.Range("E22").Value = '  +2.92%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.75%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.59'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.25%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +19.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +24.02%  '

$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.172'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.20%  '

$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.116'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +14.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '26.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.50%  '

$ws.Range("E31").Value = '  +0.13%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.93'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.97%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.13'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '51.10'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.39%  '

$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0458'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.81%  '

$ws.Range("B36").Value = 'Toncoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.74%  '

$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.06'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.64%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.23'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.41%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.62'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.86'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.117'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.76%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '122.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.08'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.82%  '

$ws.Range("E45").Value = '  +19.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.05'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.52%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.91%  '

$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.37'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.79%  '

$ws.Range("D49").Value = '2.043.90'
$ws.Range("E49").Value = '  +1.38%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0335'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.05%  '

$ws.Range("B51").Value = 'SEI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.868'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.41%  '
